$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 104.794801
$ws.Range("H2").Value = 314.384403
$ws.Range("I2").Value = 0.3872421191355361
$ws.Range("J2").Value = 0.3872421191355361
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 505.8908573333333
$ws.Range("N2").Value = 1517.672572
$ws.Range("O2").Value = 0.7018211771568338
$ws.Range("P2").Value = 0.7018211771568337
$ws.Range("Q2").Value = 53014.73172196606
$ws.Range("R2").Value = 477132.5854976946
$ws.Range("S2").Value = 0.2717747198964088
$ws.Range("T2").Value = 0.2717747198964088

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 104.794801
$ws.Range("H3").Value = 314.384403
$ws.Range("I3").Value = 0.3872421191355361
$ws.Range("J3").Value = 0.3872421191355361
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 88.00803400000001
$ws.Range("N3").Value = 264.024102
$ws.Range("O3").Value = 0.1220933352041998
$ws.Range("P3").Value = 0.1220933352041997
$ws.Range("Q3").Value = 9222.784409431235
$ws.Range("R3").Value = 83005.05968488113
$ws.Range("S3").Value = 0.04727968185679967
$ws.Range("T3").Value = 0.04727968185679966

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 104.794801
$ws.Range("H4").Value = 314.384403
$ws.Range("I4").Value = 0.3872421191355361
$ws.Range("J4").Value = 0.3872421191355361
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 126.926974
$ws.Range("N4").Value = 380.780922
$ws.Range("O4").Value = 0.1760854876389666
$ws.Range("P4").Value = 0.1760854876389665
$ws.Range("Q4").Value = 13301.28698186218
$ws.Range("R4").Value = 119711.5828367596
$ws.Range("S4").Value = 0.06818771738232766
$ws.Range("T4").Value = 0.06818771738232764

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 136.674446
$ws.Range("H5").Value = 410.023338
$ws.Range("I5").Value = 0.5050451128841343
$ws.Range("J5").Value = 0.5050451128841343
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 505.8908573333333
$ws.Range("N5").Value = 1517.672572
$ws.Range("O5").Value = 0.7018211771568338
$ws.Range("P5").Value = 0.7018211771568337
$ws.Range("Q5").Value = 69142.35266249836
$ws.Range("R5").Value = 622281.1739624853
$ws.Range("S5").Value = 0.3544513556416491
$ws.Range("T5").Value = 0.3544513556416491

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 136.674446
$ws.Range("H6").Value = 410.023338
$ws.Range("I6").Value = 0.5050451128841343
$ws.Range("J6").Value = 0.5050451128841343
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 88.00803400000001
$ws.Range("N6").Value = 264.024102
$ws.Range("O6").Value = 0.1220933352041998
$ws.Range("P6").Value = 0.1220933352041997
$ws.Range("Q6").Value = 12028.44929049916
$ws.Range("R6").Value = 108256.0436144925
$ws.Range("S6").Value = 0.06166264226060552
$ws.Range("T6").Value = 0.06166264226060551

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 136.674446
$ws.Range("H7").Value = 410.023338
$ws.Range("I7").Value = 0.5050451128841343
$ws.Range("J7").Value = 0.5050451128841343
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 126.926974
$ws.Range("N7").Value = 380.780922
$ws.Range("O7").Value = 0.1760854876389666
$ws.Range("P7").Value = 0.1760854876389665
$ws.Range("Q7").Value = 17347.67385390641
$ws.Range("R7").Value = 156129.0646851576
$ws.Range("S7").Value = 0.08893111498187971
$ws.Range("T7").Value = 0.08893111498187968

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 29.14904533333333
$ws.Range("H8").Value = 87.447136
$ws.Range("I8").Value = 0.1077127679803296
$ws.Range("J8").Value = 0.1077127679803296
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 505.8908573333333
$ws.Range("N8").Value = 1517.672572
$ws.Range("O8").Value = 0.7018211771568338
$ws.Range("P8").Value = 0.7018211771568337
$ws.Range("Q8").Value = 14746.2355341282
$ws.Range("R8").Value = 132716.1198071538
$ws.Range("S8").Value = 0.07559510161877581
$ws.Range("T8").Value = 0.07559510161877581

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 29.14904533333333
$ws.Range("H9").Value = 87.447136
$ws.Range("I9").Value = 0.1077127679803296
$ws.Range("J9").Value = 0.1077127679803296
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 88.00803400000001
$ws.Range("N9").Value = 264.024102
$ws.Range("O9").Value = 0.1220933352041998
$ws.Range("P9").Value = 0.1220933352041997
$ws.Range("Q9").Value = 2565.350172763542
$ws.Range("R9").Value = 23088.15155487187
$ws.Range("S9").Value = 0.01315101108679457
$ws.Range("T9").Value = 0.01315101108679457

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 29.14904533333333
$ws.Range("H10").Value = 87.447136
$ws.Range("I10").Value = 0.1077127679803296
$ws.Range("J10").Value = 0.1077127679803296
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 126.926974
$ws.Range("N10").Value = 380.780922
$ws.Range("O10").Value = 0.1760854876389666
$ws.Range("P10").Value = 0.1760854876389665
$ws.Range("Q10").Value = 3699.800119148822
$ws.Range("R10").Value = 33298.2010723394
$ws.Range("S10").Value = 0.01896665527475919
$ws.Range("T10").Value = 0.01896665527475919
